$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.282086133956909
$ws.Range("B1").Value = 2.933879137039185
$ws.Range("C1").Value = 5.343530178070068
$ws.Range("D1").Value = 1.853372573852539
$ws.Range("E1").Value = 1.019717574119568
